$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows above row 627; this pushes the existing rows 627-642
# down to 630-645, which already matches the target data for those rows.
$ws.Rows.Item(627).Insert()
$ws.Rows.Item(627).Insert()
$ws.Rows.Item(627).Insert()

# Fill in the constant columns (same for every row in this block) for the
# three newly inserted rows.
$constCols = @{
    1  = 10
    2  = "Vega Modelo de Temuco"
    3  = "La Araucanía"
    5  = 9
    6  = 100112028
    7  = "Sandia"
    8  = "Sin especificar"
    17 = 1
    18 = "Hortaliza"
}

foreach ($r in 627..629) {
    foreach ($col in $constCols.Keys) {
        $ws.Cells.Item($r, $col).Value = $constCols[$col]
    }
}

# Row 627
$ws.Cells.Item(627, 4).Value  = 44939          # Fecha
$ws.Cells.Item(627, 9).Value  = "Extra"        # Calidad
$ws.Cells.Item(627, 10).Value = 800            # Volumen
$ws.Cells.Item(627, 11).Value = 3500           # Precio minimo
$ws.Cells.Item(627, 12).Value = 3500           # Precio maximo
$ws.Cells.Item(627, 13).Value = 3500           # Precio promedio ponderado
$ws.Cells.Item(627, 14).Value = "$/unidad"     # Unidad de comercializacion
$ws.Cells.Item(627, 15).Value = "Región del Maule" # Origen
$ws.Cells.Item(627, 16).Value = 3500           # Precio $/Kg

# Row 628
$ws.Cells.Item(628, 4).Value  = 44939
$ws.Cells.Item(628, 9).Value  = "Primera"
$ws.Cells.Item(628, 10).Value = 2500
$ws.Cells.Item(628, 11).Value = 3000
$ws.Cells.Item(628, 12).Value = 3000
$ws.Cells.Item(628, 13).Value = 3000
$ws.Cells.Item(628, 14).Value = "$/unidad"
$ws.Cells.Item(628, 15).Value = "Región del Maule"
$ws.Cells.Item(628, 16).Value = 3000

# Row 629
$ws.Cells.Item(629, 4).Value  = 44939
$ws.Cells.Item(629, 9).Value  = "Segunda"
$ws.Cells.Item(629, 10).Value = 1000
$ws.Cells.Item(629, 11).Value = 2500
$ws.Cells.Item(629, 12).Value = 2500
$ws.Cells.Item(629, 13).Value = 2500
$ws.Cells.Item(629, 14).Value = "$/unidad"
$ws.Cells.Item(629, 15).Value = "Región del Maule"
$ws.Cells.Item(629, 16).Value = 2500
